$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 'StackedEnsemble_AllModels_4_AutoML_1_20240525_142933'

# Row 3
$ws.Range("A3").Value = 'StackedEnsemble_AllModels_6_AutoML_1_20240525_142933'

# Row 4
$ws.Range("A4").Value = 'StackedEnsemble_Best1000_1_AutoML_1_20240525_142933'
$ws.Range("B4").Value = 41454.8147379177
$ws.Range("C4").Value = 1718501664.95508
$ws.Range("D4").Value = 23294.3303148232
$ws.Range("F4").Value = 1718501664.95508

# Row 5
$ws.Range("A5").Value = 'StackedEnsemble_AllModels_3_AutoML_1_20240525_142933'

# Row 6
$ws.Range("A6").Value = 'StackedEnsemble_AllModels_2_AutoML_1_20240525_142933'

# Row 7
$ws.Range("A7").Value = 'StackedEnsemble_AllModels_1_AutoML_1_20240525_142933'

# Row 8
$ws.Range("A8").Value = 'GBM_2_AutoML_1_20240525_142933'

# Row 9
$ws.Range("A9").Value = 'StackedEnsemble_BestOfFamily_2_AutoML_1_20240525_142933'

# Row 10
$ws.Range("A10").Value = 'StackedEnsemble_BestOfFamily_4_AutoML_1_20240525_142933'
$ws.Range("B10").Value = 42046.6137318345
$ws.Range("C10").Value = 1767917726.31409
$ws.Range("D10").Value = 23824.0538579689
$ws.Range("F10").Value = 1767917726.31409

# Row 11
$ws.Range("A11").Value = 'StackedEnsemble_BestOfFamily_3_AutoML_1_20240525_142933'
$ws.Range("B11").Value = 42049.0451578926
$ws.Range("C11").Value = 1768122198.69049
$ws.Range("D11").Value = 23822.1606639521
$ws.Range("F11").Value = 1768122198.69049

# Row 12
$ws.Range("A12").Value = 'StackedEnsemble_BestOfFamily_6_AutoML_1_20240525_142933'
$ws.Range("B12").Value = 42076.909972133
$ws.Range("C12").Value = 1770466352.80299
$ws.Range("D12").Value = 23846.2304531303
$ws.Range("F12").Value = 1770466352.80299

# Row 13
$ws.Range("A13").Value = 'GBM_grid_1_AutoML_1_20240525_142933_model_16'

# Row 14
$ws.Range("A14").Value = 'GBM_grid_1_AutoML_1_20240525_142933_model_2'

# Row 15
$ws.Range("A15").Value = 'GBM_3_AutoML_1_20240525_142933'

# Row 16
$ws.Range("A16").Value = 'GBM_grid_1_AutoML_1_20240525_142933_model_1'

# Row 17
$ws.Range("A17").Value = 'GBM_grid_1_AutoML_1_20240525_142933_model_20'
$ws.Range("B17").Value = 42304.0587258997
$ws.Range("C17").Value = 1789633384.68437
$ws.Range("D17").Value = 24456.622863756
$ws.Range("E17").ClearContents()
$ws.Range("F17").Value = 1789633384.68437

# Row 18
$ws.Range("A18").Value = 'StackedEnsemble_AllModels_5_AutoML_1_20240525_142933'
$ws.Range("B18").Value = 42304.9562902227
$ws.Range("C18").Value = 1789709326.71765
$ws.Range("D18").Value = 23496.9785109574
$ws.Range("E18").Value = 3.002425906844
$ws.Range("F18").Value = 1789709326.71765

# Row 19
$ws.Range("A19").Value = 'GBM_grid_1_AutoML_1_20240525_142933_model_4'

# Row 20
$ws.Range("A20").Value = 'GBM_grid_1_AutoML_1_20240525_142933_model_17'

# Row 21
$ws.Range("A21").Value = 'GBM_1_AutoML_1_20240525_142933'

# Row 22
$ws.Range("A22").Value = 'GBM_grid_1_AutoML_1_20240525_142933_model_22'

# Row 23
$ws.Range("A23").Value = 'StackedEnsemble_BestOfFamily_1_AutoML_1_20240525_142933'

# Row 24
$ws.Range("A24").Value = 'GBM_grid_1_AutoML_1_20240525_142933_model_11'

# Row 25
$ws.Range("A25").Value = 'GBM_grid_1_AutoML_1_20240525_142933_model_7'

# Row 26
$ws.Range("A26").Value = 'GBM_5_AutoML_1_20240525_142933'

# Row 27
$ws.Range("A27").Value = 'GBM_grid_1_AutoML_1_20240525_142933_model_10'

# Row 28
$ws.Range("A28").Value = 'GBM_grid_1_AutoML_1_20240525_142933_model_24'

# Row 29
$ws.Range("A29").Value = 'GBM_grid_1_AutoML_1_20240525_142933_model_13'

# Row 30
$ws.Range("A30").Value = 'GBM_grid_1_AutoML_1_20240525_142933_model_6'

# Row 31
$ws.Range("A31").Value = 'GBM_grid_1_AutoML_1_20240525_142933_model_14'

# Row 32
$ws.Range("A32").Value = 'StackedEnsemble_BestOfFamily_5_AutoML_1_20240525_142933'
$ws.Range("B32").Value = 42808.5852150732
$ws.Range("C32").Value = 1832574968.11619
$ws.Range("D32").Value = 24166.6552363412
$ws.Range("E32").Value = 3.00193317106895
$ws.Range("F32").Value = 1832574968.11619

# Row 33
$ws.Range("A33").Value = 'GBM_grid_1_AutoML_1_20240525_142933_model_5'

# Row 34
$ws.Range("A34").Value = 'GBM_grid_1_AutoML_1_20240525_142933_model_8'

# Row 35
$ws.Range("A35").Value = 'GBM_grid_1_AutoML_1_20240525_142933_model_50'
$ws.Range("B35").Value = 42943.208398986
$ws.Range("C35").Value = 1844119147.59874
$ws.Range("D35").Value = 24869.8716158882
$ws.Range("F35").Value = 1844119147.59874

# Row 36
$ws.Range("A36").Value = 'GBM_grid_1_AutoML_1_20240525_142933_model_21'
$ws.Range("B36").Value = 42986.584473911
$ws.Range("C36").Value = 1847846444.73269
$ws.Range("D36").Value = 25293.4468413183
$ws.Range("F36").Value = 1847846444.73269

# Row 37
$ws.Range("A37").Value = 'GBM_grid_1_AutoML_1_20240525_142933_model_15'
$ws.Range("B37").Value = 42996.2334460932
$ws.Range("C37").Value = 1848676090.55095
$ws.Range("D37").Value = 25059.7531429502
$ws.Range("F37").Value = 1848676090.55095

# Row 38
$ws.Range("A38").Value = 'GBM_4_AutoML_1_20240525_142933'
$ws.Range("B38").Value = 43000.4738309286
$ws.Range("C38").Value = 1849040749.68437
$ws.Range("D38").Value = 24578.2327437266
$ws.Range("F38").Value = 1849040749.68437

# Row 39
$ws.Range("A39").Value = 'GBM_grid_1_AutoML_1_20240525_142933_model_49'
$ws.Range("B39").Value = 43782.8859403507
$ws.Range("C39").Value = 1916941101.26576
$ws.Range("D39").Value = 25513.6696901049
$ws.Range("F39").Value = 1916941101.26576

# Row 40
$ws.Range("A40").Value = 'GBM_grid_1_AutoML_1_20240525_142933_model_18'

# Row 41
$ws.Range("A41").Value = 'GBM_grid_1_AutoML_1_20240525_142933_model_19'

# Row 42
$ws.Range("A42").Value = 'DRF_1_AutoML_1_20240525_142933'

# Row 43
$ws.Range("A43").Value = 'XRT_1_AutoML_1_20240525_142933'

# Row 44
$ws.Range("A44").Value = 'GBM_grid_1_AutoML_1_20240525_142933_model_23'

# Row 45
$ws.Range("A45").Value = 'GBM_grid_1_AutoML_1_20240525_142933_model_3'

# Row 46
$ws.Range("A46").Value = 'GBM_grid_1_AutoML_1_20240525_142933_model_12'

# Row 47
$ws.Range("A47").Value = 'GBM_grid_1_AutoML_1_20240525_142933_model_9'

# Row 48
$ws.Range("A48").Value = 'DeepLearning_grid_1_AutoML_1_20240525_142933_model_1'
$ws.Range("B48").Value = 48181.4447224675
$ws.Range("C48").Value = 2321451615.54419
$ws.Range("D48").Value = 29822.1355951402
$ws.Range("F48").Value = 2321451615.54419

# Row 49
$ws.Range("A49").Value = 'DeepLearning_grid_1_AutoML_1_20240525_142933_model_3'
$ws.Range("B49").Value = 48310.5829853662
$ws.Range("C49").Value = 2333912428.38596
$ws.Range("D49").Value = 29596.305746314
$ws.Range("F49").Value = 2333912428.38596

# Row 50
$ws.Range("A50").Value = 'DeepLearning_1_AutoML_1_20240525_142933'
$ws.Range("B50").Value = 49248.6750431834
$ws.Range("C50").Value = 2425431993.50907
$ws.Range("D50").Value = 31150.5162145537
$ws.Range("F50").Value = 2425431993.50907

# Row 51
$ws.Range("A51").Value = 'DeepLearning_grid_1_AutoML_1_20240525_142933_model_4'
$ws.Range("B51").Value = 49364.875213671
$ws.Range("C51").Value = 2436890904.86131
$ws.Range("D51").Value = 31326.4734378545
$ws.Range("E51").ClearContents()
$ws.Range("F51").Value = 2436890904.86131

# Row 52
$ws.Range("A52").Value = 'GBM_grid_1_AutoML_1_20240525_142933_model_51'
$ws.Range("B52").Value = 49547.0324099614
$ws.Range("C52").Value = 2454908420.63377
$ws.Range("D52").Value = 31956.4374772204
$ws.Range("E52").Value = 3.32808142069402
$ws.Range("F52").Value = 2454908420.63377

# Row 53
$ws.Range("A53").Value = 'DeepLearning_grid_3_AutoML_1_20240525_142933_model_3'
$ws.Range("B53").Value = 49547.6627480308
$ws.Range("C53").Value = 2454970883.7926
$ws.Range("D53").Value = 32115.6869891946
$ws.Range("E53").Value = 3.26901564396188
$ws.Range("F53").Value = 2454970883.7926

# Row 54
$ws.Range("A54").Value = 'DeepLearning_grid_2_AutoML_1_20240525_142933_model_3'
$ws.Range("B54").Value = 49600.1662468616
$ws.Range("C54").Value = 2460176491.71631
$ws.Range("D54").Value = 32200.2549882723
$ws.Range("E54").Value = 3.25126344450862
$ws.Range("F54").Value = 2460176491.71631

# Row 55
$ws.Range("A55").Value = 'DeepLearning_grid_1_AutoML_1_20240525_142933_model_9'
$ws.Range("B55").Value = 51178.1784954011
$ws.Range("C55").Value = 2619205954.10713
$ws.Range("D55").Value = 32348.6244768845
$ws.Range("F55").Value = 2619205954.10713

# Row 56
$ws.Range("A56").Value = 'DeepLearning_grid_1_AutoML_1_20240525_142933_model_2'
$ws.Range("B56").Value = 53643.635548867
$ws.Range("C56").Value = 2877639634.89967
$ws.Range("D56").Value = 33917.388149328
$ws.Range("E56").ClearContents()
$ws.Range("F56").Value = 2877639634.89967

# Row 57
$ws.Range("A57").Value = 'DeepLearning_grid_3_AutoML_1_20240525_142933_model_2'
$ws.Range("B57").Value = 54407.399043952
$ws.Range("C57").Value = 2960165070.72783
$ws.Range("D57").Value = 39296.08939981
$ws.Range("E57").Value = 3.50364695768501
$ws.Range("F57").Value = 2960165070.72783

# Row 58
$ws.Range("A58").Value = 'DeepLearning_grid_2_AutoML_1_20240525_142933_model_2'
$ws.Range("B58").Value = 55202.7299580206
$ws.Range("C58").Value = 3047341394.81815
$ws.Range("D58").Value = 40137.7739927866
$ws.Range("E58").Value = 3.5102735467676
$ws.Range("F58").Value = 3047341394.81815

# Row 59
$ws.Range("A59").Value = 'DeepLearning_grid_3_AutoML_1_20240525_142933_model_1'
$ws.Range("B59").Value = 56368.672361003
$ws.Range("C59").Value = 3177427223.7421
$ws.Range("D59").Value = 44342.1504942428
$ws.Range("E59").Value = 3.58504348729328
$ws.Range("F59").Value = 3177427223.7421

# Row 60
$ws.Range("A60").Value = 'GLM_1_AutoML_1_20240525_142933'

# Row 61
$ws.Range("A61").Value = 'DeepLearning_grid_2_AutoML_1_20240525_142933_model_1'
$ws.Range("B61").Value = 61945.7377602217
$ws.Range("C61").Value = 3837274426.65815
$ws.Range("D61").Value = 51122.2305832154
$ws.Range("E61").Value = 3.6626927483468
$ws.Range("F61").Value = 3837274426.65815
